# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 114
    $ws.Range("F3").Value = 21
    $ws.Range("F4").Value = 973
    $ws.Range("F5").Value = 490
}
